$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Unity 6.0 Addressable Scene loading change -> add a new "SceneData" table
# block (rows 119-123) to the design spec sheet, mirroring the existing
# "Projectile" table block layout (rows 109-114/117-118).
# ---------------------------------------------------------------------------

# Row 119: blank spacer row (same shape as rows 111 / 116 / 121)
$ws.Range("A119").Formula = "=ROW()"
$ws.Range("A119").VerticalAlignment = -4108

# Row 120: new table-header row ("G_1024" / SceneData), formatted like row 110
$ws.Range("A110:K110").Copy()
$ws.Range("A120:K120").PasteSpecial(-4122)
$ws.Range("A120").Formula = "=ROW()"
$ws.Range("A120").VerticalAlignment = -4108
$ws.Range("B120").Value = "G_1024"
$ws.Range("C120").Value = "Y"
$ws.Range("D120").Value = "Y"
$ws.Range("E120").Value = "씬 데이터"
$ws.Range("F120").Value = "SceneData"

# Row 121: blank spacer row
$ws.Range("A121").Formula = "=ROW()"
$ws.Range("A121").VerticalAlignment = -4108

# Row 122: new "SceneName" PK column row, formatted like row 112
$ws.Range("A112:I112").Copy()
$ws.Range("A122:I122").PasteSpecial(-4122)
$ws.Range("K112").Copy()
$ws.Range("K122").PasteSpecial(-4122)
$ws.Range("A122").Formula = "=ROW()"
$ws.Range("A122").VerticalAlignment = -4108
$ws.Range("C122").Value = "Y"
$ws.Range("D122").Value = "Y"
$ws.Range("E122").Value = "SceneName"
$ws.Range("F122").Value = "SceneName"
$ws.Range("G122").Value = "string"
$ws.Range("H122").Value = "PK"
$ws.Range("K122").Value = "-"

# Row 123: new "Path" column row, formatted like row 118
$ws.Range("A118:G118").Copy()
$ws.Range("A123:G123").PasteSpecial(-4122)
$ws.Range("I118").Copy()
$ws.Range("I123").PasteSpecial(-4122)
$ws.Range("A123").Formula = "=ROW()"
$ws.Range("A123").VerticalAlignment = -4108
$ws.Range("C123").Value = "Y"
$ws.Range("D123").Value = "Y"
$ws.Range("E123").Value = "Path"
$ws.Range("F123").Value = "Path"
$ws.Range("G123").Value = "string"

# Extend the list-validation range on column G to cover the new rows, same
# as it already does for G117:G118
$existing = $ws.Range("G52:G53,G57:G67,G40:G48,G71:G76,G80:G85,G89:G91,G95:G100,G104:G108,G15:G36,G112:G114,G5:G11,G117:G118,G122:G123")
$existing.Validation.Delete()
$existing.Validation.Add(3, 1, 1, "=DataType")
$existing.Validation.InCellDropdown = $true
$existing.Validation.IgnoreBlank = $true

# Keep the selection anchored near the new rows, matching the authored edit
$ws.Application.Goto($ws.Range("A91"), $true)
$ws.Range("J117:K118").Select()
